$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Partially revert reduced Food weight (column F, "粮") for rows 21-23.
$ws.Range("F21").Value = 0.7
$ws.Range("F22").Value = 0.7
$ws.Range("F23").Value = 0.7

# Update the visible selection / scroll position to match the saved view.
$ws.Range("F24").Select()
$excel.ActiveWindow.ScrollRow = 4
